$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 2 ("H") ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 191
$wsOff.Range("C2").Value = 130
$wsOff.Range("D2").Value = 49
$wsOff.Range("E2").Value = 27
$wsOff.Range("G2").Value = 3

# --- DEF sheet: row 2 ("H") ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 248
$wsDef.Range("C2").Value = 168
$wsDef.Range("D2").Value = 59
$wsDef.Range("E2").Value = 32
$wsDef.Range("F2").Value = 5
